$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Files-tab Cypher query in B4: drop the `File Type` and `Breed`
# output columns (ICDC Breed 1-14 script correction).
$newQuery = "MATCH (f:file)-->(parent)`nWITH DISTINCT f, parent`nMATCH (f)-[*]->(c:case)<--(demo:demographic)`nWHERE demo.breed IN ['Australian Shepherd']`nOPTIONAL MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)`nOPTIONAL MATCH (samp:sample)-->(c)`nWITH DISTINCT f, parent, c, demo, diag, s`nRETURN  coalesce(f.file_name, '') AS ``File Name``,`n        coalesce(labels(parent)[0], '') AS ``Association``,`n        coalesce(f.file_description, '') AS ``Description``,`n        coalesce(f.file_format, '') AS ``Format``,`n        coalesce(f.file_size, '') AS ``Size``,`n        coalesce(c.case_id, '') AS ``Case ID``,`n        coalesce(diag.disease_term,'') AS Diagnosis , `n        coalesce(s.clinical_study_designation,'') AS ``Study Code``"

$ws.Range("B4").Value = $newQuery

# Move the active selection to D4 (was D12).
$ws.Range("D4").Select()

# The shorter query text re-wraps to fewer lines, so the row shrinks.
$ws.Rows("4:4").RowHeight = 217.5
